$wb = $excel.ActiveWorkbook

# Rename the "WasteAllocationx" sheet to "WasteAllocation"
$wasteSheet = $wb.Worksheets.Item("WasteAllocationx")
$wasteSheet.Name = "WasteAllocation"

# Move the active tab from "Processes" (previously selected) to "WasteAllocation",
# and update its selection to A2 (was D9).
$wasteSheet.Activate() | Out-Null
$wasteSheet.Range("A2").Select() | Out-Null
